$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1975806451612903
$ws.Range("C2").Value = 0.5685483870967742
$ws.Range("J2").Value = 0.008064516129032258
$ws.Range("P2").Value = 0.1491935483870968
$ws.Range("S2").Value = 0.07661290322580645
$ws.Range("B3").Value = 0.00684931506849315
$ws.Range("C3").Value = 0.0410958904109589
$ws.Range("J3").Value = 0.0410958904109589
$ws.Range("P3").Value = 0.7397260273972602
$ws.Range("S3").Value = 0.1712328767123288
$ws.Range("J4").Value = 0.05128205128205128
$ws.Range("P4").Value = 0.6410256410256411
$ws.Range("S4").Value = 0.3076923076923077
$ws.Range("B6").Value = 0.04672897196261682
$ws.Range("D6").Value = 0.01869158878504673
$ws.Range("E6").Value = 0.004672897196261682
$ws.Range("F6").Value = 0.0514018691588785
$ws.Range("J6").Value = 0.2102803738317757
$ws.Range("O6").Value = 0.01869158878504673
$ws.Range("Q6").Value = 0.1682242990654206
$ws.Range("R6").Value = 0.1121495327102804
$ws.Range("S6").Value = 0.3691588785046729
$ws.Range("B7").Value = 0.07741935483870968
$ws.Range("D7").Value = 0.01290322580645161
$ws.Range("F7").Value = 0.03870967741935484
$ws.Range("J7").Value = 0.09677419354838709
$ws.Range("Q7").Value = 0.1548387096774194
$ws.Range("R7").Value = 0.1032258064516129
$ws.Range("S7").Value = 0.5161290322580645
$ws.Range("B8").Value = 0.07943925233644859
$ws.Range("D8").Value = 0.02336448598130841
$ws.Range("E8").Value = 0.004672897196261682
$ws.Range("F8").Value = 0.08644859813084112
$ws.Range("J8").Value = 0.102803738317757
$ws.Range("O8").Value = 0.01168224299065421
$ws.Range("Q8").Value = 0.1869158878504673
$ws.Range("R8").Value = 0.1004672897196262
$ws.Range("S8").Value = 0.4042056074766355
$ws.Range("B9").Value = 0.07784431137724551
$ws.Range("D9").Value = 0.01197604790419162
$ws.Range("F9").Value = 0.08982035928143713
$ws.Range("J9").Value = 0.05389221556886228
$ws.Range("O9").Value = 0.02395209580838323
$ws.Range("Q9").Value = 0.2215568862275449
$ws.Range("R9").Value = 0.1197604790419162
$ws.Range("S9").Value = 0.4011976047904192
$ws.Range("B10").Value = 0.1017087062652563
$ws.Range("D10").Value = 0.01708706265256306
$ws.Range("E10").Value = 0.002441008950366151
$ws.Range("F10").Value = 0.06916192026037429
$ws.Range("J10").Value = 0.09601301871440195
$ws.Range("O10").Value = 0.01464605370219691
$ws.Range("Q10").Value = 0.2131814483319772
$ws.Range("R10").Value = 0.1082180634662327
$ws.Range("S10").Value = 0.3775427176566314
$ws.Range("G11").Value = 0.1862348178137652
$ws.Range("J11").Value = 0.1052631578947368
$ws.Range("K11").Value = 0.2388663967611336
$ws.Range("L11").Value = 0.4615384615384616
$ws.Range("S11").Value = 0.008097165991902834
$ws.Range("G12").Value = 0.8205128205128205
$ws.Range("J12").Value = 0.1196581196581197
$ws.Range("K12").Value = 0.008547008547008548
$ws.Range("L12").Value = 0.02564102564102564
$ws.Range("S12").Value = 0.02564102564102564
$ws.Range("G13").Value = 0.5600000000000001
$ws.Range("J13").Value = 0.44
$ws.Range("F15").Value = 0.02564102564102564
$ws.Range("H15").Value = 0.1435897435897436
$ws.Range("I15").Value = 0.1128205128205128
$ws.Range("J15").Value = 0.3538461538461539
$ws.Range("K15").Value = 0.05641025641025641
$ws.Range("M15").Value = 0.01025641025641026
$ws.Range("O15").Value = 0.07179487179487179
$ws.Range("S15").Value = 0.2256410256410256
$ws.Range("F16").Value = 0.01197604790419162
$ws.Range("H16").Value = 0.1437125748502994
$ws.Range("I16").Value = 0.07784431137724551
$ws.Range("J16").Value = 0.4131736526946108
$ws.Range("K16").Value = 0.1437125748502994
$ws.Range("M16").Value = 0.005988023952095809
$ws.Range("N16").Value = 0.005988023952095809
$ws.Range("O16").Value = 0.05988023952095808
$ws.Range("S16").Value = 0.1377245508982036
$ws.Range("F17").Value = 0.00684931506849315
$ws.Range("H17").Value = 0.1986301369863014
$ws.Range("I17").Value = 0.07534246575342465
$ws.Range("J17").Value = 0.4726027397260274
$ws.Range("K17").Value = 0.07990867579908675
$ws.Range("M17").Value = 0.0136986301369863
$ws.Range("N17").Value = 0.00228310502283105
$ws.Range("O17").Value = 0.0593607305936073
$ws.Range("S17").Value = 0.091324200913242
$ws.Range("F18").Value = 0.008510638297872341
$ws.Range("H18").Value = 0.1574468085106383
$ws.Range("I18").Value = 0.06808510638297872
$ws.Range("J18").Value = 0.4936170212765957
$ws.Range("K18").Value = 0.09787234042553192
$ws.Range("M18").Value = 0.01276595744680851
$ws.Range("N18").Value = 0.00425531914893617
$ws.Range("O18").Value = 0.0425531914893617
$ws.Range("S18").Value = 0.1148936170212766
$ws.Range("F19").Value = 0.01306620209059233
$ws.Range("H19").Value = 0.2212543554006969
$ws.Range("I19").Value = 0.07317073170731707
$ws.Range("J19").Value = 0.4268292682926829
$ws.Range("K19").Value = 0.07752613240418119
$ws.Range("M19").Value = 0.009581881533101045
$ws.Range("N19").Value = 0.0008710801393728223
$ws.Range("O19").Value = 0.07578397212543554
$ws.Range("S19").Value = 0.1019163763066202
